$d = $word.ActiveDocument

# Pass 1: old text -> unique placeholder tokens (avoids collisions between
# an old value and another pair's new value)
$d.Content.Find.Execute("2024-10-07 Monday", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN0@@", 2) | Out-Null
$d.Content.Find.Execute("81-24=57", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN1@@", 2) | Out-Null
$d.Content.Find.Execute("60-35=25", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN2@@", 2) | Out-Null
$d.Content.Find.Execute("54+18=72", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN3@@", 2) | Out-Null
$d.Content.Find.Execute("40-3=37", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN4@@", 2) | Out-Null
$d.Content.Find.Execute("91-42=49", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN5@@", 2) | Out-Null
$d.Content.Find.Execute("16+27=43", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN6@@", 2) | Out-Null
$d.Content.Find.Execute("93-86=7", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN7@@", 2) | Out-Null
$d.Content.Find.Execute("35-28=7", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN8@@", 2) | Out-Null
$d.Content.Find.Execute("36+19=55", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN9@@", 2) | Out-Null
$d.Content.Find.Execute("53-4=49", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN10@@", 2) | Out-Null
$d.Content.Find.Execute("70-27=43", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN11@@", 2) | Out-Null
$d.Content.Find.Execute("8+26=34", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN12@@", 2) | Out-Null
$d.Content.Find.Execute("74-69=5", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN13@@", 2) | Out-Null
$d.Content.Find.Execute("88+6=94", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN14@@", 2) | Out-Null
$d.Content.Find.Execute("72-34=38", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN15@@", 2) | Out-Null
$d.Content.Find.Execute("6+35=41", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN16@@", 2) | Out-Null
$d.Content.Find.Execute("72-7=65", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN17@@", 2) | Out-Null
$d.Content.Find.Execute("56-49=7", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN18@@", 2) | Out-Null
$d.Content.Find.Execute("25+57=82", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN19@@", 2) | Out-Null
$d.Content.Find.Execute("91-45=46", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN20@@", 2) | Out-Null
$d.Content.Find.Execute("55-36=19", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN21@@", 2) | Out-Null
$d.Content.Find.Execute("96-39=57", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN22@@", 2) | Out-Null
$d.Content.Find.Execute("84-77=7", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN23@@", 2) | Out-Null
$d.Content.Find.Execute("56+35=91", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN24@@", 2) | Out-Null
$d.Content.Find.Execute("76+19=95", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN25@@", 2) | Out-Null
$d.Content.Find.Execute("71-45=26", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN26@@", 2) | Out-Null
$d.Content.Find.Execute("47-8=39", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN27@@", 2) | Out-Null
$d.Content.Find.Execute("90-61=29", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN28@@", 2) | Out-Null
$d.Content.Find.Execute("18+27=45", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN29@@", 2) | Out-Null
$d.Content.Find.Execute("46+49=95", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN30@@", 2) | Out-Null
$d.Content.Find.Execute("32-4=28", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN31@@", 2) | Out-Null
$d.Content.Find.Execute("84-67=17", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN32@@", 2) | Out-Null
$d.Content.Find.Execute("64+28=92", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN33@@", 2) | Out-Null
$d.Content.Find.Execute("18+4=22", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN34@@", 2) | Out-Null
$d.Content.Find.Execute("85-58=27", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN35@@", 2) | Out-Null
$d.Content.Find.Execute("6+78=84", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN36@@", 2) | Out-Null
$d.Content.Find.Execute("78+4=82", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN37@@", 2) | Out-Null
$d.Content.Find.Execute("18+47=65", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN38@@", 2) | Out-Null
$d.Content.Find.Execute("84-28=56", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN39@@", 2) | Out-Null
$d.Content.Find.Execute("80-41=39", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN40@@", 2) | Out-Null
$d.Content.Find.Execute("33-9=24", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN41@@", 2) | Out-Null
$d.Content.Find.Execute("44+29=73", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN42@@", 2) | Out-Null
$d.Content.Find.Execute("52-34=18", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN43@@", 2) | Out-Null
$d.Content.Find.Execute("90-83=7", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN44@@", 2) | Out-Null
$d.Content.Find.Execute("26-18=8", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN45@@", 2) | Out-Null
$d.Content.Find.Execute("94-45=49", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN46@@", 2) | Out-Null
$d.Content.Find.Execute("14+49=63", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN47@@", 2) | Out-Null
$d.Content.Find.Execute("48+35=83", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN48@@", 2) | Out-Null
$d.Content.Find.Execute("29+69=98", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN49@@", 2) | Out-Null
$d.Content.Find.Execute("64-9=55", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN50@@", 2) | Out-Null
$d.Content.Find.Execute("95-29=66", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN51@@", 2) | Out-Null
$d.Content.Find.Execute("60-33=27", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN52@@", 2) | Out-Null
$d.Content.Find.Execute("81-27=54", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN53@@", 2) | Out-Null
$d.Content.Find.Execute("37+16=53", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN54@@", 2) | Out-Null
$d.Content.Find.Execute("36-29=7", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN55@@", 2) | Out-Null
$d.Content.Find.Execute("93-66=27", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN56@@", 2) | Out-Null
$d.Content.Find.Execute("25+58=83", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN57@@", 2) | Out-Null
$d.Content.Find.Execute("94-47=47", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN58@@", 2) | Out-Null
$d.Content.Find.Execute("3+88=91", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN59@@", 2) | Out-Null
$d.Content.Find.Execute("66+16=82", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN60@@", 2) | Out-Null
$d.Content.Find.Execute("39+59=98", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN61@@", 2) | Out-Null
$d.Content.Find.Execute("39+3=42", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN62@@", 2) | Out-Null
$d.Content.Find.Execute("4+59=63", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN63@@", 2) | Out-Null
$d.Content.Find.Execute("61-7=54", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN64@@", 2) | Out-Null
$d.Content.Find.Execute("62-59=3", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN65@@", 2) | Out-Null
$d.Content.Find.Execute("26-9=17", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN66@@", 2) | Out-Null
$d.Content.Find.Execute("91-27=64", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN67@@", 2) | Out-Null
$d.Content.Find.Execute("56-37=19", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN68@@", 2) | Out-Null
$d.Content.Find.Execute("38+29=67", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN69@@", 2) | Out-Null
$d.Content.Find.Execute("41-34=7", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN70@@", 2) | Out-Null
$d.Content.Find.Execute("62-19=43", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN71@@", 2) | Out-Null
$d.Content.Find.Execute("29+63=92", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN72@@", 2) | Out-Null
$d.Content.Find.Execute("85-49=36", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN73@@", 2) | Out-Null
$d.Content.Find.Execute("70-21=49", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN74@@", 2) | Out-Null
$d.Content.Find.Execute("38+59=97", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN75@@", 2) | Out-Null
$d.Content.Find.Execute("67+29=96", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN76@@", 2) | Out-Null
$d.Content.Find.Execute("8+54=62", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN77@@", 2) | Out-Null
$d.Content.Find.Execute("97-58=39", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN78@@", 2) | Out-Null
$d.Content.Find.Execute("37-8=29", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN79@@", 2) | Out-Null
$d.Content.Find.Execute("9+76=85", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN80@@", 2) | Out-Null
$d.Content.Find.Execute("41-33=8", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN81@@", 2) | Out-Null
$d.Content.Find.Execute("19+35=54", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN82@@", 2) | Out-Null
$d.Content.Find.Execute("94-88=6", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN83@@", 2) | Out-Null
$d.Content.Find.Execute("71-23=48", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN84@@", 2) | Out-Null
$d.Content.Find.Execute("71-65=6", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN85@@", 2) | Out-Null
$d.Content.Find.Execute("55-29=26", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN86@@", 2) | Out-Null
$d.Content.Find.Execute("81-49=32", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN87@@", 2) | Out-Null
$d.Content.Find.Execute("80-9=71", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN88@@", 2) | Out-Null
$d.Content.Find.Execute("95-76=19", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN89@@", 2) | Out-Null
$d.Content.Find.Execute("34-26=8", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN90@@", 2) | Out-Null
$d.Content.Find.Execute("29+9=38", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN91@@", 2) | Out-Null
$d.Content.Find.Execute("25+59=84", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN92@@", 2) | Out-Null
$d.Content.Find.Execute("54+29=83", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN93@@", 2) | Out-Null
$d.Content.Find.Execute("34+49=83", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN94@@", 2) | Out-Null
$d.Content.Find.Execute("69+3=72", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN95@@", 2) | Out-Null
$d.Content.Find.Execute("45+38=83", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN96@@", 2) | Out-Null
$d.Content.Find.Execute("56-17=39", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN97@@", 2) | Out-Null
$d.Content.Find.Execute("45-29=16", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN98@@", 2) | Out-Null
$d.Content.Find.Execute("17+14=31", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN99@@", 2) | Out-Null
$d.Content.Find.Execute("65+16=81", $true, $true, $false, $false, $false, $true, 1, $false, "@@TOKEN100@@", 2) | Out-Null

# Pass 2: placeholder tokens -> final new text
$d.Content.Find.Execute("@@TOKEN0@@", $true, $true, $false, $false, $false, $true, 1, $false, "2024-10-08 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN1@@", $true, $true, $false, $false, $false, $true, 1, $false, "7+46=53", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN2@@", $true, $true, $false, $false, $false, $true, 1, $false, "19+4=23", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN3@@", $true, $true, $false, $false, $false, $true, 1, $false, "74-19=55", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN4@@", $true, $true, $false, $false, $false, $true, 1, $false, "73-66=7", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN5@@", $true, $true, $false, $false, $false, $true, 1, $false, "96-68=28", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN6@@", $true, $true, $false, $false, $false, $true, 1, $false, "58+27=85", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN7@@", $true, $true, $false, $false, $false, $true, 1, $false, "97-18=79", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN8@@", $true, $true, $false, $false, $false, $true, 1, $false, "92-35=57", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN9@@", $true, $true, $false, $false, $false, $true, 1, $false, "41-12=29", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN10@@", $true, $true, $false, $false, $false, $true, 1, $false, "41-39=2", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN11@@", $true, $true, $false, $false, $false, $true, 1, $false, "34+39=73", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN12@@", $true, $true, $false, $false, $false, $true, 1, $false, "48+46=94", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN13@@", $true, $true, $false, $false, $false, $true, 1, $false, "6+26=32", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN14@@", $true, $true, $false, $false, $false, $true, 1, $false, "90-31=59", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN15@@", $true, $true, $false, $false, $false, $true, 1, $false, "18+63=81", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN16@@", $true, $true, $false, $false, $false, $true, 1, $false, "88+3=91", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN17@@", $true, $true, $false, $false, $false, $true, 1, $false, "94-29=65", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN18@@", $true, $true, $false, $false, $false, $true, 1, $false, "95-56=39", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN19@@", $true, $true, $false, $false, $false, $true, 1, $false, "55-46=9", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN20@@", $true, $true, $false, $false, $false, $true, 1, $false, "33+38=71", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN21@@", $true, $true, $false, $false, $false, $true, 1, $false, "58-9=49", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN22@@", $true, $true, $false, $false, $false, $true, 1, $false, "24+7=31", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN23@@", $true, $true, $false, $false, $false, $true, 1, $false, "65-57=8", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN24@@", $true, $true, $false, $false, $false, $true, 1, $false, "37+7=44", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN25@@", $true, $true, $false, $false, $false, $true, 1, $false, "74-58=16", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN26@@", $true, $true, $false, $false, $false, $true, 1, $false, "55-16=39", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN27@@", $true, $true, $false, $false, $false, $true, 1, $false, "44+49=93", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN28@@", $true, $true, $false, $false, $false, $true, 1, $false, "92-25=67", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN29@@", $true, $true, $false, $false, $false, $true, 1, $false, "92-39=53", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN30@@", $true, $true, $false, $false, $false, $true, 1, $false, "7+44=51", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN31@@", $true, $true, $false, $false, $false, $true, 1, $false, "59+12=71", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN32@@", $true, $true, $false, $false, $false, $true, 1, $false, "55-26=29", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN33@@", $true, $true, $false, $false, $false, $true, 1, $false, "49+2=51", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN34@@", $true, $true, $false, $false, $false, $true, 1, $false, "37+58=95", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN35@@", $true, $true, $false, $false, $false, $true, 1, $false, "34+39=73", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN36@@", $true, $true, $false, $false, $false, $true, 1, $false, "11-3=8", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN37@@", $true, $true, $false, $false, $false, $true, 1, $false, "24-6=18", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN38@@", $true, $true, $false, $false, $false, $true, 1, $false, "40-21=19", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN39@@", $true, $true, $false, $false, $false, $true, 1, $false, "59+25=84", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN40@@", $true, $true, $false, $false, $false, $true, 1, $false, "59+6=65", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN41@@", $true, $true, $false, $false, $false, $true, 1, $false, "25-19=6", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN42@@", $true, $true, $false, $false, $false, $true, 1, $false, "24+39=63", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN43@@", $true, $true, $false, $false, $false, $true, 1, $false, "47+29=76", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN44@@", $true, $true, $false, $false, $false, $true, 1, $false, "9+14=23", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN45@@", $true, $true, $false, $false, $false, $true, 1, $false, "53-16=37", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN46@@", $true, $true, $false, $false, $false, $true, 1, $false, "70-44=26", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN47@@", $true, $true, $false, $false, $false, $true, 1, $false, "38+26=64", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN48@@", $true, $true, $false, $false, $false, $true, 1, $false, "49+4=53", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN49@@", $true, $true, $false, $false, $false, $true, 1, $false, "33-4=29", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN50@@", $true, $true, $false, $false, $false, $true, 1, $false, "6+26=32", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN51@@", $true, $true, $false, $false, $false, $true, 1, $false, "86+5=91", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN52@@", $true, $true, $false, $false, $false, $true, 1, $false, "93-18=75", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN53@@", $true, $true, $false, $false, $false, $true, 1, $false, "4+39=43", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN54@@", $true, $true, $false, $false, $false, $true, 1, $false, "87+8=95", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN55@@", $true, $true, $false, $false, $false, $true, 1, $false, "68+4=72", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN56@@", $true, $true, $false, $false, $false, $true, 1, $false, "54+27=81", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN57@@", $true, $true, $false, $false, $false, $true, 1, $false, "27+65=92", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN58@@", $true, $true, $false, $false, $false, $true, 1, $false, "19+45=64", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN59@@", $true, $true, $false, $false, $false, $true, 1, $false, "63+8=71", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN60@@", $true, $true, $false, $false, $false, $true, 1, $false, "14+79=93", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN61@@", $true, $true, $false, $false, $false, $true, 1, $false, "90-43=47", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN62@@", $true, $true, $false, $false, $false, $true, 1, $false, "26+36=62", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN63@@", $true, $true, $false, $false, $false, $true, 1, $false, "8+59=67", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN64@@", $true, $true, $false, $false, $false, $true, 1, $false, "7+45=52", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN65@@", $true, $true, $false, $false, $false, $true, 1, $false, "30-2=28", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN66@@", $true, $true, $false, $false, $false, $true, 1, $false, "34-26=8", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN67@@", $true, $true, $false, $false, $false, $true, 1, $false, "56-47=9", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN68@@", $true, $true, $false, $false, $false, $true, 1, $false, "66+8=74", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN69@@", $true, $true, $false, $false, $false, $true, 1, $false, "14+48=62", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN70@@", $true, $true, $false, $false, $false, $true, 1, $false, "86+7=93", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN71@@", $true, $true, $false, $false, $false, $true, 1, $false, "26+8=34", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN72@@", $true, $true, $false, $false, $false, $true, 1, $false, "23+38=61", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN73@@", $true, $true, $false, $false, $false, $true, 1, $false, "7+87=94", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN74@@", $true, $true, $false, $false, $false, $true, 1, $false, "8+83=91", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN75@@", $true, $true, $false, $false, $false, $true, 1, $false, "62-8=54", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN76@@", $true, $true, $false, $false, $false, $true, 1, $false, "28+8=36", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN77@@", $true, $true, $false, $false, $false, $true, 1, $false, "85-77=8", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN78@@", $true, $true, $false, $false, $false, $true, 1, $false, "7+66=73", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN79@@", $true, $true, $false, $false, $false, $true, 1, $false, "92-37=55", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN80@@", $true, $true, $false, $false, $false, $true, 1, $false, "56+18=74", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN81@@", $true, $true, $false, $false, $false, $true, 1, $false, "91-45=46", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN82@@", $true, $true, $false, $false, $false, $true, 1, $false, "9+25=34", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN83@@", $true, $true, $false, $false, $false, $true, 1, $false, "65-26=39", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN84@@", $true, $true, $false, $false, $false, $true, 1, $false, "59+13=72", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN85@@", $true, $true, $false, $false, $false, $true, 1, $false, "27+27=54", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN86@@", $true, $true, $false, $false, $false, $true, 1, $false, "64-46=18", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN87@@", $true, $true, $false, $false, $false, $true, 1, $false, "46-7=39", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN88@@", $true, $true, $false, $false, $false, $true, 1, $false, "76-28=48", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN89@@", $true, $true, $false, $false, $false, $true, 1, $false, "84-45=39", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN90@@", $true, $true, $false, $false, $false, $true, 1, $false, "10-5=5", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN91@@", $true, $true, $false, $false, $false, $true, 1, $false, "61-28=33", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN92@@", $true, $true, $false, $false, $false, $true, 1, $false, "19+4=23", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN93@@", $true, $true, $false, $false, $false, $true, 1, $false, "29+2=31", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN94@@", $true, $true, $false, $false, $false, $true, 1, $false, "80-1=79", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN95@@", $true, $true, $false, $false, $false, $true, 1, $false, "36-27=9", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN96@@", $true, $true, $false, $false, $false, $true, 1, $false, "79+7=86", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN97@@", $true, $true, $false, $false, $false, $true, 1, $false, "90-55=35", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN98@@", $true, $true, $false, $false, $false, $true, 1, $false, "15+19=34", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN99@@", $true, $true, $false, $false, $false, $true, 1, $false, "50-43=7", 2) | Out-Null
$d.Content.Find.Execute("@@TOKEN100@@", $true, $true, $false, $false, $false, $true, 1, $false, "41-35=6", 2) | Out-Null
